# Remove the two oldest date rows (2025-09-21 and 2025-09-22) from the
# "Chart" sheet. Excel shifts the remaining rows up automatically, and the
# now-unused shared strings ("2025-09-21", the blank placeholder, and
# "2025-09-22") are dropped from the shared-string table on save.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Chart")

$ws.Rows.Item(2).Delete()
$ws.Rows.Item(2).Delete()
